$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.232.01'
$ws.Range("E2").Value = '  +4.13%  '
$ws.Range("D3").Value = '3.621.60'
$ws.Range("E3").Value = '  +6.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.73%  '
$ws.Range("D7").Value = '3.615.60'
$ws.Range("E7").Value = '  +6.90%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.608'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.93%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  +5.04%  '
$ws.Range("E11").Value = '  +2.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '50.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.35%  '
$ws.Range("E13").Value = '  +2.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '705.63'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.89%  '
$ws.Range("D15").Value = '4.205.54'
$ws.Range("E15").Value = '  +6.99%  '
$ws.Range("E16").Value = '  +3.69%  '
$ws.Range("D17").Value = '72.215.91'
$ws.Range("E17").Value = '  +3.98%  '
$ws.Range("D18").Value = '3.603.44'
$ws.Range("E18").Value = '  +6.18%  '
$ws.Range("E19").Value = '  +1.77%  '
$ws.Range("E20").Value = '  +4.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.89%  '
$ws.Range("E22").Value = '  +3.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '105.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.66%  '
$ws.Range("E27").Value = '  +4.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +16.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '593.66'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.49%  '
$ws.Range("E34").Value = '  +2.12%  '
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.77'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.00%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.145'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.41%  '
$ws.Range("D39").Value = '3.650.55'
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0782'
$ws.Range("E40").Value = '  +7.97%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.06'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.80'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.91%  '
$ws.Range("E44").Value = '  +6.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.349'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.93%  '
$ws.Range("E47").Value = '  +4.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.49'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.24%  '
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '133.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.03%  '
